$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-9). Sending cluster / Ligand / Receptor / Target cluster
# plus the 16 numeric metric columns (E..T), replacing the previous TPM-derived values.
$data = @(
    @("FAPs",  "Angpt4", "Tek", "ECs",           3, 1, 1.091893666666667,  3.275681,            0.8167354825840876, 0.8167354825840876, 3, 1, 56.16842399999999,          168.505272,  0.8021403195141494,  0.8021403195141494,  61.329946432248,    551.9695178902319, 0.6551364609585431,  0.6551364609585431),
    @("FAPs",  "Angpt4", "Tek", "FAPs",          3, 1, 1.091893666666667,  3.275681,            0.8167354825840876, 0.8167354825840876, 3, 1, 10.582537,                   31.747611,   0.1511290331103168,  0.1511290331103167,  11.55500512756567,  103.995046148091,  0.1234324437898211,  0.1234324437898211),
    @("FAPs",  "Angpt4", "Tek", "MuSCs",         3, 1, 1.091893666666667,  3.275681,            0.8167354825840876, 0.8167354825840876, 3, 1, 3.262797333333333,          9.788392,    0.04659595390231912, 0.04659595390231912, 3.562627743883556,  32.063649694952,   0.03805656889687651, 0.03805656889687651),
    @("FAPs",  "Angpt4", "Tek", "Resolving-Mac", 3, 1, 1.091893666666667,  3.275681,            0.8167354825840876, 0.8167354825840876, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.01029837709944444, 0.09268539389500001, 0.0001100089388468628, 0.0001100089388468628),
    @("MuSCs", "Angpt4", "Tek", "ECs",           3, 1, 0.2450063333333334, 0.7350190000000001, 0.1832645174159124, 0.1832645174159125, 3, 1, 56.16842399999999,          168.505272,  0.8021403195141494,  0.8021403195141494,  13.761619613352,    123.854576520168,  0.1470038585556064,  0.1470038585556064),
    @("MuSCs", "Angpt4", "Tek", "FAPs",          3, 1, 0.2450063333333334, 0.7350190000000001, 0.1832645174159124, 0.1832645174159125, 3, 1, 10.582537,                   31.747611,   0.1511290331103168,  0.1511290331103167,  2.592788587734334,  23.335097289609,   0.02769658932049565, 0.02769658932049565),
    @("MuSCs", "Angpt4", "Tek", "MuSCs",         3, 1, 0.2450063333333334, 0.7350190000000001, 0.1832645174159124, 0.1832645174159125, 3, 1, 3.262797333333333,          9.788392,    0.04659595390231912, 0.04659595390231912, 0.7994060110497778, 7.194654099448001, 0.008539385005442617, 0.008539385005442617),
    @("MuSCs", "Angpt4", "Tek", "Resolving-Mac", 3, 1, 0.2450063333333334, 0.7350190000000001, 0.1832645174159124, 0.1832645174159125, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.002310818067222222, 0.020797362605, 0.00002468453436774894, 0.00002468453436774895)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# The old sheet had 9 data rows (row 10); the new data only has 8 (rows 2-9).
$ws.Rows.Item(10).Delete()
